# Update the ParkingGarage worksheet to the latest "legacy version found in UW
# servers": rename a few garages, swap a couple of icon file names, and add
# two brand-new rows (Cedar Garage East / West).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "South Campus Garage" -> "S1 Garage" (icon SCG.png -> S1.png)
$ws.Range("A4").Value = "S1 Garage"
$ws.Range("D4").Value = "S1.png"

# Row 5: Padelford Garage keeps its name but gets a new icon file name
$ws.Range("D5").Value = "PDL.png"

# Row 8: 4545 Garage keeps its name but gets a new icon file name
$ws.Range("D8").Value = "G4545.png"

# Row 9: "Tower Garage B" -> "Tower Garage B (W45)" (icon TGB.png -> W45.png)
$ws.Range("A9").Value = "Tower Garage B (W45)"
$ws.Range("D9").Value = "W45.png"

# Row 10: "Tower Garage A" -> "Tower Garage A (W46)" (icon TGA.png -> W46.png)
$ws.Range("A10").Value = "Tower Garage A (W46)"
$ws.Range("D10").Value = "W46.png"

# New row 11: Cedar Garage East
$ws.Range("A11").Value = "Cedar Garage East"
$ws.Range("B11").Value = 47.656922074768502
$ws.Range("C11").Value = -122.315618991851
$ws.Range("D11").Value = "CGE.png"

# New row 12: Cedar Garage West
$ws.Range("A12").Value = "Cedar Garage West"
$ws.Range("B12").Value = 47.657446001010499
$ws.Range("C12").Value = -122.316579222679
$ws.Range("D12").Value = "CGW.png"

# Column A widens slightly to fit the new, longer garage names.
$ws.Columns.Item(1).ColumnWidth = 19.5

# Update the worksheet selection to match the latest save.
$ws.Range("C18").Select()
